# Remove the paragraph "Occurrence: Role instance." entirely (including its
# paragraph mark), leaving the preceding "Metaclass: Role." paragraph and the
# following (empty) paragraph untouched.
$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($text -eq "Occurrence: Role instance.") {
        $p.Range.Delete()
        break
    }
}
